# Applies the "updated pricing data" edit to the pricingData table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# ------------------------------------------------------------------
# 1. Correct a handful of existing data points (2019-2021 rows).
# ------------------------------------------------------------------

# Row 21 - Overleaf / 2019
$ws.Cells.Item(21, 7).Value = 1        # G21 managementFeatures
$ws.Cells.Item(21, 11).Value = 3       # K21 numberOfPlans

# Row 26 - Clockify / 2022
$ws.Cells.Item(26, 3).Value = 4        # C26 informationFeatures
$ws.Cells.Item(26, 4).Value = 5        # D26 integrationFeatures
$ws.Cells.Item(26, 5).Value = 23       # E26 toolFeatures
$ws.Cells.Item(26, 6).Value = 4        # F26 automationFeatures
$ws.Cells.Item(26, 7).Value = 7        # G26 managementFeatures
$ws.Cells.Item(26, 11).Value = 5       # K26 numberOfPlans

# Row 28 - Salescloud / 2019
$ws.Cells.Item(28, 5).Value = 40       # E28 toolFeatures
$ws.Cells.Item(28, 6).Value = 2        # F28 automationFeatures
$ws.Cells.Item(28, 7).Value = 3        # G28 managementFeatures

# Row 31 - Jira / 2020
$ws.Cells.Item(31, 3).Value = 5        # C31 informationFeatures
$ws.Cells.Item(31, 4).Value = 3        # D31 integrationFeatures
$ws.Cells.Item(31, 5).Value = 33       # E31 toolFeatures
$ws.Cells.Item(31, 7).Value = 1        # G31 managementFeatures
$ws.Cells.Item(31, 10).Value = 22      # J31 numberOfCommonFeatures

# Row 34 - Overleaf / 2020
$ws.Cells.Item(34, 3).Value = 1        # C34 informationFeatures
$ws.Cells.Item(34, 5).Value = 6        # E34 toolFeatures
$ws.Cells.Item(34, 11).Value = 3       # K34 numberOfPlans

# Row 36 - Postman / 2020
$ws.Cells.Item(36, 3).Value = 2        # C36 informationFeatures
$ws.Cells.Item(36, 4).Value = 3        # D36 integrationFeatures
$ws.Cells.Item(36, 5).Value = 8        # E36 toolFeatures

# Row 38 - Salescloud / 2020
$ws.Cells.Item(38, 9).Value = 11       # I38 numberOfAddOns

# Row 47 - Overleaf / 2021
$ws.Cells.Item(47, 3).Value = 1        # C47 informationFeatures
$ws.Cells.Item(47, 5).Value = 7        # E47 toolFeatures
$ws.Cells.Item(47, 11).Value = 4       # K47 numberOfPlans

# ------------------------------------------------------------------
# 2. Append three new 2022 rows (Github, Jira, Overleaf) to the table.
# ------------------------------------------------------------------

function Add-PricingRow([string]$name, [int]$year, [int]$info, [int]$integ, [int]$tool, [int]$auto, [int]$mgmt, [int]$addons, [int]$common, [int]$plans) {
    $lastRow = $lo.Range.Row + $lo.Range.Rows.Count - 1
    $srcRange = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 11))
    $srcRange.Copy()
    $lo.ListRows.Add() | Out-Null
    $newRow = $lo.Range.Row + $lo.Range.Rows.Count - 1
    $dstRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 11))
    $dstRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Rows.Item($newRow).RowHeight = 18

    $ws.Cells.Item($newRow, 1).Value = $name
    $ws.Cells.Item($newRow, 2).Value = $year
    $ws.Cells.Item($newRow, 3).Value = $info
    $ws.Cells.Item($newRow, 4).Value = $integ
    $ws.Cells.Item($newRow, 5).Value = $tool
    $ws.Cells.Item($newRow, 6).Value = $auto
    $ws.Cells.Item($newRow, 7).Value = $mgmt
    $ws.Cells.Item($newRow, 8).Formula = "=SUM(C" + $newRow + ":G" + $newRow + ")"
    $ws.Cells.Item($newRow, 9).Value = $addons
    $ws.Cells.Item($newRow, 10).Value = $common
    $ws.Cells.Item($newRow, 11).Value = $plans
}

Add-PricingRow "Github"   2022 3 5 21 9 9 0 17 3
Add-PricingRow "Jira"     2022 3 1 19 1 7 1 12 3
Add-PricingRow "Overleaf" 2022 1 5 8  0 0 0 3  4

# ------------------------------------------------------------------
# 3. Filter the table down to year 2023 (hides every other row).
# ------------------------------------------------------------------

$lo.Range.AutoFilter(2, @("2023"), 7) | Out-Null

# ------------------------------------------------------------------
# 4. Scroll the sheet view down one row (below the frozen header).
# ------------------------------------------------------------------

$ws.Activate()
$ws.Range("A2").Select()
